# Edit: draft day 13 blog updated
# - Correct Day 13 (row 14) Ayats range from "114 - 129" to "114 - 130"
# - Add new Day 14 (row 15) entry: date, ayats, content, author, tags
# - Update selection / row height to reflect the new entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 (Day 14): new entry (ayats range set first so it claims the
#     next shared-string slot ahead of Day 13's corrected range) ---
$ws.Range("C15").Value = "Surah Baqarah, 131 - 144"

# --- Row 14 (Day 13): fix the Ayats range ---
# (Content + Tags text is unchanged; only the shared-string slot shifts because
#  the old "114 - 129" string becomes orphaned once replaced below.)
$ws.Range("C14").Value = "Surah Baqarah, 114 - 130"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 43844
$ws.Range("B15").NumberFormat = $ws.Range("B14").NumberFormat

$ws.Range("E15").Value = "Qasim Ali"

$ws.Range("F15").Value = "Doubting Allah (swt), Colors of Islam, Prophet (PBUH) a blessing"

$d15Content = @"
h1: Colors of Allah (swt)
p: A human being has minute cells in body. Each cell is an independent industry that powers on its own. Creation of one cell will need a factory as big as New York city. Still making it come alive is not possible. Two feet distance between me and my computer has infinite particles. In air alone which I can not see. And what to talk of seeing, I can only see the Visible Spectrum of Electromagnetic spectrum. <a href="https://earthsky.org/space/what-is-the-electromagnetic-spectrum"> Rest 99 % is not even visible to us.</a> The thing that can see is an Eye. How does one cell manufactured inside human body gets to decide it will be converted to a glass like structure. How does it get connected to retina (back of the eye). How does all of this happen? Coincidence? Evolution after a Big Bang? Initially it was only lava and nothing. And then it got converted to nature. 
h3: Who is Allah (swt)?
quote: Have you not considered how Allah presents an example, [making] a good word like a good tree, whose root is firmly fixed and its branches [high] in the sky? <br> - Surah Ibrahim Verse 24
p: I have been looking for the right person to follow. Seriously, I felt Steve Jobs should be followed being such a great presenter. May be Allama Iqbal. Anthony Robbins. Someone who is both vocal and practical. Someone following whom gives my life a meaning. But who should be followed? Who really changed the face of this planet?
h3: Who should I follow?
p: A person, who rises in ruins. Abused, degraded, humiliated in public. Who does that. Runs from a society of chaos. People are drinking. Brothels are open. He runs away for 2 years. When he is a tradesman, his enemies respect him. When he is a leader, his people put fear into their hearts. When he is bullied, he laughs. When he is with kids, he is a kid. How does one does all this. So who should any body follow. Any sane person in this world should follow 1 guy who does all this. Because it will make anybody a success in this life. So when everything he does, is meant to be followed. 
p: <b>Why can’t I follow the One he talked about his entire life? </b> Just because I suddenly became more logical and practical? Sanity leads to Allah (swt).
h3: Doubts are fine. Ibrahim (as) also had doubts. 
quote: And [remember] when Ibrahim said, “My Lord! Show me how You give life to the dead.” Allah said, “What, do you not believe [in resurrection]?” He answered, “Certainly [I believe, I am asking this] so that my heart may be at ease.” <br> Surah Baqarah Verse 260
p: So in such times, the only solution is to hold tight to the rope. There is no other way of succeeding in this life. 
h3: Why Prophet Muhammad (saw) is a blessing for Muslims?
p: Allah (swt) says:-
quote: We have certainly seen the turning of your face, [O Muhammad], toward the heaven, and We will surely turn you to a qiblah with which you will be pleased. So turn your face toward Al-Masjid Al-Haram. And wherever you [believers] are, turn your faces toward it [in prayer]. Indeed, those who have been given the Scripture well know that it is the truth from their Lord. And Allah is not unaware of what they do. <br> Surah Baqarah, verse 144
p.b-left: Prophet Mohammad (peace be upon him) turned towards the new Qiblah without hesitation, accurately facing the Kaaba without the use of any scientific instrument or compass. This took place in the second year after Hijra, the migration. Most interpretations date it to the middle of the month of Shaban.
p.b-left: <a href=“http://saudigazette.com.sa/article/125309”>Link here</a>
p: Prophet PBUH did not wait for somebody to bring compass to calculate accurate direction of Kaba before turning to Makkah. Neither did he wait for the prayer to finish. He just acted upon what was instructed. He has established a path that can not be doubted. The real blessing is his sunnah which sets forth the colorful path. 
p: This path tells us to embrace the calamities with an open heart. 
"@
$ws.Range("D15").Value = $d15Content

$ws.Rows.Item(15).RowHeight = 409.6

# --- Update active selection to the new row ---
$ws.Range("D15").Select()

Write-Host "Edit applied"
